$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so numeric-looking
# strings like "1.00" or "0.161" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "68.185.07"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "3.759.45"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "595.20"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "167.63"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.756.52"
$ws.Range("E7").Value = "  -1.20%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "36.74"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.384.30"
$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.748.31"
$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "68.080.76"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "18.24"
$ws.Range("E18").Value = "  -3.11%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "7.07"
$ws.Range("E19").Value = "  -3.63%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.112"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "10.86"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "470.35"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.704"
$ws.Range("E23").Value = "  -4.42%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "83.27"
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("E25").Value = "  -5.01%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  -0.65%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.900.44"
$ws.Range("E30").Value = "  -1.43%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  -4.03%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.41"
$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "30.10"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "9.30"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  --%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.712.08"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -9.08%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.138"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "5.82"
$ws.Range("E42").Value = "  -1.94%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.309"
$ws.Range("E44").Value = "  -2.22%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.94"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "8.60"
$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "400.83"
$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "45.30"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "144.87"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "25.91"
$ws.Range("E51").Value = "  +2.26%  "
